$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.143.48"
$ws.Range("E2").Value = "  -3.31%  "
$ws.Range("D3").Value = "1.926.47"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.39"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4726"
$ws.Range("E7").Value = "  -5.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4069"
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.08"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08442"
$ws.Range("E10").Value = "  -8.79%  "
$ws.Range("E11").Value = "  -4.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.25"
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("D13").Value = "1.927.82"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.516"
$ws.Range("E14").Value = "  -5.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.098"
$ws.Range("E15").Value = "  -5.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.66"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001065"
$ws.Range("E18").Value = "  -3.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06582"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.17"
$ws.Range("E20").Value = "  -5.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.767"
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("D23").Value = "28.136.08"
$ws.Range("E23").Value = "  -3.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.42"
$ws.Range("E24").Value = "  -4.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.276"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").Value = "2.139.64"
$ws.Range("E26").Value = "  -3.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.35"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.13"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.155"
$ws.Range("E29").Value = "  -4.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.735"
$ws.Range("E30").Value = "  -9.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.83"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9735"
$ws.Range("E32").Value = "  -7.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09596"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.440"
$ws.Range("E34").Value = "  -5.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.557"
$ws.Range("E35").Value = "  -4.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.638"
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.079"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02313"
$ws.Range("E38").Value = "  -5.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06172"
$ws.Range("E39").Value = "  -3.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.236"
$ws.Range("E40").Value = "  -6.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6174"
$ws.Range("E41").Value = "  -4.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.05"
$ws.Range("E42").Value = "  -3.85%  "
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1904"
$ws.Range("E44").Value = "  -4.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.301"
$ws.Range("E45").Value = "  -5.61%  "
$ws.Range("E46").Value = "  -5.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.79"
$ws.Range("E47").Value = "  -4.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.036"
$ws.Range("E48").Value = "  -7.30%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06814"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.83"
$ws.Range("E51").Value = "  -3.13%  "
